$d = $word.ActiveDocument

function Do-Split($findText, $replaceText, $label) {
    $d.Content.Find.ClearFormatting()
    $ok = $d.Content.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Replacement failed for: $label"
    } else {
        Write-Output "OK: $label"
    }
}

$find0 = "1.Generic biotechnological process: schematic representation; description of the main stages.2.Equipment sterilization: terminology; sterilization by physical agents; sterilization by chemical agents.3.Media sterilization by steam heating: kinetics of thermal destruction of microorganisms; destruction of nutrient media; calculation of sterilization time for batch processes; design of sterilization systems for continuous processes.4.Sterilization by filtration: microbial aerosols; air samplers; filter sizing; media sterilization.5.Kinetics and stoichiometry of microbial growth and product formation: transformation rates and conversion factors; classification of fermentative processes based on cell growth and product formation rates; influence of substrate concentration on cell growth rate; stoichiometry of microbial growth and product formation."
$replace0 = "1.Generic biotechnological process: schematic representation; description of the main stages.^l2.Equipment sterilization: terminology; sterilization by physical agents; sterilization by chemical agents.^l3.Media sterilization by steam heating: kinetics of thermal destruction of microorganisms; destruction of nutrient media; calculation of sterilization time for batch processes; design of sterilization systems for continuous processes.^l4.Sterilization by filtration: microbial aerosols; air samplers; filter sizing; media sterilization.^l5.Kinetics and stoichiometry of microbial growth and product formation: transformation rates and conversion factors; classification of fermentative processes based on cell growth and product formation rates; influence of substrate concentration on cell growth rate; stoichiometry of microbial growth and product formation."
Do-Split $find0 $replace0 "Programa (EN, italic)"

$find1 = "A nota final (NF) será composta pelas médias M1  e M2,calculadas conforme segue:M1=P1+a1×T1M2=P2+a2×T2Em que:-P1 e P2 são as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).-T1 e T2 são as notas médias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.-a1 e a2 são os fatores multiplicadores das notas médias dos trabalhos, a serem definidos pelo docente antes do início de cada turma com base nas atividades específicas a serem propostas. Os valores serão ≥0,1, sendo informados aos alunos no início do semestre. Em todos os casos, os valores máximos para M1 e M2 serão “dez”, sendo desconsideradas pontuações superiores.O cálculo de NF será feito conforme segue:NF=(M1+2×M2)/3Serão aprovados os alunos que obtiverem NF maior ou igual 5,0."
$replace1 = "A nota final (NF) será composta pelas médias M1  e M2,calculadas conforme segue:^lM1=P1+a1×T1^lM2=P2+a2×T2^lEm que:^l-P1 e P2 são as notas da primeira e da segunda prova escrita, respectivamente (notas de zero a dez).^l-T1 e T2 são as notas médias dos trabalhos (notas de zero a dez) realizados antes da primeira e da segunda prova escrita, respectivamente.^l-a1 e a2 são os fatores multiplicadores das notas médias dos trabalhos, a serem definidos pelo docente antes do início de cada turma com base nas atividades específicas a serem propostas. Os valores serão ≥0,1, sendo informados aos alunos no início do semestre. ^lEm todos os casos, os valores máximos para M1 e M2 serão “dez”, sendo desconsideradas pontuações superiores.^lO cálculo de NF será feito conforme segue:^lNF=(M1+2×M2)/3^lSerão aprovados os alunos que obtiverem NF maior ou igual 5,0."
Do-Split $find1 $replace1 "Criterio (NF formula)"

$find2 = "Será oferecido um programa de recuperação, sendo este avaliado por uma prova escrita final (PR). A média de recuperação (MR) será calculada conforme segue: MR=(NF+PR)/2Serão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
$replace2 = "Será oferecido um programa de recuperação, sendo este avaliado por uma prova escrita final (PR). A média de recuperação (MR) será calculada conforme segue: ^lMR=(NF+PR)/2^lSerão aprovados os alunos que obtiverem MR maior ou igual a 5,0."
Do-Split $find2 $replace2 "Norma de recuperacao (MR formula)"

$find3 = "ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. Biotecnologia Industrial. Volume 1: Fundamentos. 2ª Edição. São Paulo: Blucher, 2020. ISBN 978-85-212-1897-5 (e-Book); 978-85-212-1898-2 (Impresso).ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. 2ª Edição. São Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).BORZANI, W. Processo Biotecnológico Industrial Genérico. In: BORZANI, W.; SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E. Biotecnologia Industrial. Volume 1: Fundamentos. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0278-3.DORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914SHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706"
$replace3 = "ALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. Biotecnologia Industrial. Volume 1: Fundamentos. 2ª Edição. São Paulo: Blucher, 2020. ISBN 978-85-212-1897-5 (e-Book); 978-85-212-1898-2 (Impresso).^lALTERTHUM, F.; SCHMIDELL, W.; LIMA, U. A.; MORAES. M. O. (Org.). Biotecnologia Industrial. Volume 2: Engenharia Bioquímica. 2ª Edição. São Paulo: Blucher, 2021. p. 37-52.  ISBN 978-65-5506-019-5 (e-Book); 978-65-5506-018-8 (Impresso).^lBORZANI, W. Processo Biotecnológico Industrial Genérico. In: BORZANI, W.; SCHMIDELL, W.; LIMA, U. A.; AQUARONE, E. Biotecnologia Industrial. Volume 1: Fundamentos. São Paulo: Editora Edgard Blücher Ltda, 2001. ISBN 978-85-212-0278-3.^lDORAN P.M.; MORRISSEY, K.; CARLSON, R. P. Bioprocess Engineering Principles, 3rd edition, Academic Press, 2024. ISBN 978-0128221914^lSHULER, M. L.; KARGI, F.; DELISA, M. Bioprocess Engineering: Basic Concepts (3rd Edition) (Prentice Hall International Series in the Physical and Chemical Engineering Sciences) 3rd Edition. Prentice Hall; 3 edition, 2017. ISBN: 978-0137062706"
Do-Split $find3 $replace3 "Bibliografia"
